$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new trailing columns (AD:AF) holding the team's season record,
# mirroring the header style already used by the other header cells.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-43) gets the same season record values.
for ($row = 2; $row -le 43; $row++) {
    $ws.Range("AD$row").Value = 85
    $ws.Range("AE$row").Value = 77
    $ws.Range("AF$row").Value = 0
}
